# New crime data collected
# Updates the weekly CompStat report (72nd Precinct) for the week of
# 6/10/2024 - 6/16/2024 (Volume 31, Number 24), with refreshed crime
# complaint statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text: "Volume 31   Number  23" -> "... Number  24"
# and the report date range 6/3/2024-6/9/2024 -> 6/10/2024-6/16/2024
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  24"
$ws.Range("C9").Value = "Report Covering the Week  6/10/2024  Through  6/16/2024"

# ---------------------------------------------------------------------
# Column E (Week to Date %Chg) no longer needs the extra width it had
# ---------------------------------------------------------------------
$ws.Columns("E").ColumnWidth = 6.168446

# ---------------------------------------------------------------------
# Cells that flip between "no complaints" (text placeholders "0" /
# "***.*") and real numeric figures need their number format switched
# along with the value, so copy the format from a cell that already
# carries the desired style, then (re)assign the value.
# ---------------------------------------------------------------------

# C14: was numeric 1 -> now text "0" (no murders this week)
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "0"
$ws.Range("D14").Copy() | Out-Null
$ws.Range("C14").PasteSpecial(-4122) | Out-Null

# C18: was text "0" -> now numeric 4
$ws.Range("D18").Copy() | Out-Null
$ws.Range("C18").PasteSpecial(-4122) | Out-Null
$ws.Range("C18").Value = 4

# C22: was text "0" -> now numeric 2
$ws.Range("F22").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122) | Out-Null
$ws.Range("C22").Value = 2

# D33: was numeric 1 -> now text "0"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0"
$ws.Range("C33").Copy() | Out-Null
$ws.Range("D33").PasteSpecial(-4122) | Out-Null

# E33: was numeric -100 -> now text "***.*"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "***.*"
$ws.Range("C33").Copy() | Out-Null
$ws.Range("E33").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Straightforward numeric updates (counts + computed %-change columns)
# ---------------------------------------------------------------------

# Row 14
$ws.Range("M14").Value = -80
$ws.Range("N14").Value = -88.888888888888

# Row 15
$ws.Range("F15").Value = 2
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 12
$ws.Range("K15").Value = 20
$ws.Range("L15").Value = 9.090909090909
$ws.Range("M15").Value = 50
$ws.Range("N15").Value = -25

# Row 16
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 21
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = 162.5
$ws.Range("I16").Value = 83
$ws.Range("J16").Value = 63
$ws.Range("K16").Value = 31.746031746031
$ws.Range("L16").Value = 13.698630136986
$ws.Range("M16").Value = 3.75
$ws.Range("N16").Value = -83.531746031746

# Row 17
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 15
$ws.Range("E17").Value = -60
$ws.Range("F17").Value = 21
$ws.Range("G17").Value = 34
$ws.Range("H17").Value = -38.235294117647
$ws.Range("I17").Value = 146
$ws.Range("J17").Value = 123
$ws.Range("K17").Value = 18.699186991869
$ws.Range("L17").Value = 1.388888888888
$ws.Range("M17").Value = 78.048780487804
$ws.Range("N17").Value = -38.912133891213

# Row 18
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 300
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = 11.111111111111
$ws.Range("I18").Value = 81
$ws.Range("J18").Value = 60
$ws.Range("K18").Value = 35
$ws.Range("L18").Value = -28.947368421052
$ws.Range("M18").Value = -30.769230769230
$ws.Range("N18").Value = -84.055118110236

# Row 19
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = 225
$ws.Range("F19").Value = 38
$ws.Range("G19").Value = 34
$ws.Range("H19").Value = 11.764705882352
$ws.Range("I19").Value = 230
$ws.Range("J19").Value = 259
$ws.Range("K19").Value = -11.196911196911
$ws.Range("L19").Value = -43.209876543209
$ws.Range("M19").Value = 46.496815286624
$ws.Range("N19").Value = 15

# Row 20
$ws.Range("C20").Value = 3
$ws.Range("E20").Value = -40
$ws.Range("G20").Value = 16
$ws.Range("H20").Value = -31.25
$ws.Range("I20").Value = 74
$ws.Range("J20").Value = 76
$ws.Range("K20").Value = -2.631578947368
$ws.Range("L20").Value = 12.121212121212
$ws.Range("M20").Value = 27.586206896551
$ws.Range("N20").Value = -84.188034188034

# Row 21
$ws.Range("C21").Value = 30
$ws.Range("D21").Value = 28
$ws.Range("E21").Value = 7.142857142857
$ws.Range("F21").Value = 104
$ws.Range("G21").Value = 102
$ws.Range("H21").Value = 1.960784313725
$ws.Range("I21").Value = 627
$ws.Range("J21").Value = 592
$ws.Range("K21").Value = 5.912162162162
$ws.Range("L21").Value = -22.878228782287
$ws.Range("M21").Value = 23.668639053254
$ws.Range("N21").Value = -67.746913580246

# Row 22 (remaining numeric siblings of C22, handled above)
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 8
$ws.Range("K22").Value = -38.461538461538
$ws.Range("L22").Value = -57.894736842105
$ws.Range("M22").Value = -50

# Row 24
$ws.Range("C24").Value = 31
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = 29.166666666666
$ws.Range("F24").Value = 90
$ws.Range("G24").Value = 95
$ws.Range("H24").Value = -5.263157894736
$ws.Range("I24").Value = 449
$ws.Range("J24").Value = 534
$ws.Range("K24").Value = -15.917602996254
$ws.Range("L24").Value = -29.067930489731
$ws.Range("M24").Value = 32.058823529411

# Row 25
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 13
$ws.Range("E25").Value = -30.769230769230
$ws.Range("F25").Value = 41
$ws.Range("H25").Value = -16.326530612244
$ws.Range("I25").Value = 208
$ws.Range("J25").Value = 293
$ws.Range("K25").Value = -29.010238907849
$ws.Range("L25").Value = -43.783783783783

# Row 26
$ws.Range("D26").Value = 13
$ws.Range("E26").Value = -38.461538461538
$ws.Range("F26").Value = 35
$ws.Range("G26").Value = 48
$ws.Range("H26").Value = -27.083333333333
$ws.Range("I26").Value = 226
$ws.Range("J26").Value = 220
$ws.Range("K26").Value = 2.727272727272
$ws.Range("L26").Value = 24.861878453038
$ws.Range("M26").Value = -23.389830508474

# Row 27
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = -33.333333333333
$ws.Range("I27").Value = 17
$ws.Range("K27").Value = 13.333333333333

# Row 28
$ws.Range("C28").Value = 2
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 6
$ws.Range("G28").Value = 10
$ws.Range("H28").Value = -40
$ws.Range("I28").Value = 28
$ws.Range("J28").Value = 40
$ws.Range("K28").Value = -30
$ws.Range("L28").Value = -9.677419354838

# Row 33 (D33/E33 text conversion handled above)
$ws.Range("I33").Value = 4
$ws.Range("K33").Value = 300
$ws.Range("L33").Value = 100

Write-Output "Applied weekly crime data refresh."
